# Auto-generated edit script: update cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.912.19"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.872.57"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5088"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07194"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8941"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.72"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.889.80"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07519"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.235"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008510"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "26.987.49"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.018"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "2.131.65"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.394"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.788"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.085"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.719"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.721"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09146"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05069"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7487"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.969"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.156"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.234"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.527"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5631"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02000"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.615"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.547"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1478"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4742"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.570"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.92"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.11"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.15%  "
